# Fruta / hortaliza, semanal
# Insert 3 new weekly-report rows for Femacal de La Calera - Chirimoya
# right before the existing row 177, pushing the remaining data (old
# rows 177-268) down to rows 180-271.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 177 (shifts 177-268 -> 180-271)
$ws.Range("A177:A179").EntireRow.Insert()

# --- New row 177: Chirimoya, Especial, Provincia del Elquí ---
$ws.Range("A177").Value = 3
$ws.Range("B177").Value = "Femacal de La Calera"
$ws.Range("C177").Value = "Coquimbo"
$ws.Range("D177").Value = 44879
$ws.Range("E177").Value = 5
$ws.Range("F177").Value = "Fruta"
$ws.Range("G177").Value = 100107
$ws.Range("H177").Value = "Otros"
$ws.Range("I177").Value = 100107002
$ws.Range("J177").Value = "Chirimoya"
$ws.Range("K177").Value = "Cultivar IV Región"
$ws.Range("L177").Value = "Especial"
$ws.Range("M177").Value = 56
$ws.Range("N177").Value = 26000
$ws.Range("O177").Value = 26000
$ws.Range("P177").Value = 26000
$ws.Range("Q177").Value = "$/bandeja 10 kilos"
$ws.Range("R177").Value = "Provincia del Elquí"
$ws.Range("S177").Value = 2600
$ws.Range("T177").Value = 10

# --- New row 178: Chirimoya, Primera, Provincia del Elquí ---
$ws.Range("A178").Value = 3
$ws.Range("B178").Value = "Femacal de La Calera"
$ws.Range("C178").Value = "Coquimbo"
$ws.Range("D178").Value = 44879
$ws.Range("E178").Value = 5
$ws.Range("F178").Value = "Fruta"
$ws.Range("G178").Value = 100107
$ws.Range("H178").Value = "Otros"
$ws.Range("I178").Value = 100107002
$ws.Range("J178").Value = "Chirimoya"
$ws.Range("K178").Value = "Cultivar IV Región"
$ws.Range("L178").Value = "Primera"
$ws.Range("M178").Value = 56
$ws.Range("N178").Value = 23000
$ws.Range("O178").Value = 23000
$ws.Range("P178").Value = 23000
$ws.Range("Q178").Value = "$/bandeja 10 kilos"
$ws.Range("R178").Value = "Provincia del Elquí"
$ws.Range("S178").Value = 2300
$ws.Range("T178").Value = 10

# --- New row 179: Chirimoya, Segunda, Provincia del Elquí ---
$ws.Range("A179").Value = 3
$ws.Range("B179").Value = "Femacal de La Calera"
$ws.Range("C179").Value = "Coquimbo"
$ws.Range("D179").Value = 44879
$ws.Range("E179").Value = 5
$ws.Range("F179").Value = "Fruta"
$ws.Range("G179").Value = 100107
$ws.Range("H179").Value = "Otros"
$ws.Range("I179").Value = 100107002
$ws.Range("J179").Value = "Chirimoya"
$ws.Range("K179").Value = "Cultivar IV Región"
$ws.Range("L179").Value = "Segunda"
$ws.Range("M179").Value = 60
$ws.Range("N179").Value = 20000
$ws.Range("O179").Value = 20000
$ws.Range("P179").Value = 20000
$ws.Range("Q179").Value = "$/bandeja 10 kilos"
$ws.Range("R179").Value = "Provincia del Elquí"
$ws.Range("S179").Value = 2000
$ws.Range("T179").Value = 10
